# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worksheet holds an "Estado de Cuenta" (account statement) table.
# Row 16 held the single existing debtor record (LINDA MARISOL MONTOYA
# GARCIA). The update:
#   - keeps the debtor-table row (row 16) but now shows a DIFFERENT
#     worker (SKARLYS DE JESUS NUNEZ GARCIA) with an updated "Salario
#     Basico" value,
#   - inserts 13 new detail rows right after it for two more workers
#     (GINA PAOLA SERRANO PRADA over 6 periods, JOSE GREGORIO CASTRO
#     MARTINEZ over the same 6 periods) and relocates the original
#     LINDA record to the bottom of the new table (row 29),
#   - refreshes the summary counters (Valor Mora total, worker count,
#     period count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 13 new rows right after the existing data row
#    (row 16). Everything from the old row 17 onward (the blank
#    spacer rows + the signature footer) shifts down by 13 rows.
# ---------------------------------------------------------------------
$ws.Rows("17:29").Insert()

# Copy row 16's formatting (fonts/borders/number formats) into the
# freshly inserted rows so every new record looks like the existing
# table row.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J29").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Row 16 now shows the new worker SKARLYS DE JESUS NUNEZ GARCIA
#    (same document type/period, updated "Salario Basico").
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45487196"
$ws.Range("D16").Value = "SKARLYS DE JESUS NUÑEZ GARCIA"
$ws.Range("E16").Value = "2403"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 126000

# ---------------------------------------------------------------------
# 3. New rows 17-22: GINA PAOLA SERRANO PRADA, six periods.
# ---------------------------------------------------------------------
$ginaPeriods = @("2507","2506","2505","2504","2503","2502")
$r = 17
foreach ($p in $ginaPeriods) {
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1128055034"
    $ws.Range("D$r").Value = "GINA PAOLA SERRANO PRADA"
    $ws.Range("E$r").Value = $p
    $ws.Range("F$r").Value = 138041
    $ws.Range("G$r").Value = 3451020
    $r++
}

# ---------------------------------------------------------------------
# 4. New rows 23-28: JOSE GREGORIO CASTRO MARTINEZ, same six periods.
# ---------------------------------------------------------------------
$josePeriods = @("2507","2506","2505","2504","2503","2502")
$r = 23
foreach ($p in $josePeriods) {
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "9144427"
    $ws.Range("D$r").Value = "JOSE GREGORIO CASTRO MARTINEZ"
    $ws.Range("E$r").Value = $p
    $ws.Range("F$r").Value = 46400
    $ws.Range("G$r").Value = 1160000
    $r++
}

# ---------------------------------------------------------------------
# 5. Row 29: the record that used to live in row 16 (LINDA MARISOL
#    MONTOYA GARCIA) is preserved, just relocated to the end of the
#    table with its original values untouched.
# ---------------------------------------------------------------------
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1047407520"
$ws.Range("D29").Value = "LINDA MARISOL MONTOYA GARCIA"
$ws.Range("E29").Value = "2403"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1432260

# ---------------------------------------------------------------------
# 6. Refresh the summary counters above the table.
#    Valor Mora total = sum of the "Valor Mora" (F) column.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1210646
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 7

# ---------------------------------------------------------------------
# 7. Column D ("Nombre Trabajador") now holds some longer names
#    (JOSE GREGORIO CASTRO MARTINEZ / SKARLYS DE JESUS NUNEZ GARCIA)
#    than before, so re-fit its width.
# ---------------------------------------------------------------------
$ws.Columns("D:D").AutoFit()

Write-Output "Estados de cuenta actualizados"
